$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data: add the "lifting events" (NPI lifted) entries ---
$ws.Range("A51").Value = 2
$ws.Range("A61").Value = 2
$ws.Range("A76").Value = 2

# --- Theme: background/light-1 color changed from a pale green tint back to pure white ---
$wb.Theme.ThemeColorScheme.Colors(2).RGB = 16777215

# --- View state: selection moved to E71 (closest controllable proxy for the author's view change) ---
$ws.Range("E71").Select()
